$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cut the UNIT column (D) and re-insert it after QTY (before the old PRICE
# column, which is column I before this move). This mirrors "Insert Cut
# Cells" in the Excel UI: Style/Color/Size/Qty all shift one column left,
# and Unit now sits just before Price/Surcharge.
$ws.Columns("D:D").Cut()
$ws.Columns("I:I").Insert()

# Insert a new blank column right after STYLE (column D) that will hold
# the DESCRIPTION header.
$ws.Columns("E:E").Insert()

# Re-apply the column widths to match the final layout (the shifted/new
# columns don't automatically keep the right widths after the cut/insert).
$ws.Columns("A").ColumnWidth = 9.592447916666666
$ws.Columns("B").ColumnWidth = 21.451822916666668
$ws.Columns("C").ColumnWidth = 20.307291666666668
$ws.Columns("D").ColumnWidth = 11.451822916666666
$ws.Columns("E").ColumnWidth = 11.166666666666666
$ws.Columns("F").ColumnWidth = 16.736979166666668
$ws.Columns("G").ColumnWidth = 19.166666666666668
$ws.Columns("H").ColumnWidth = 10.451822916666666
$ws.Columns("I").ColumnWidth = 10.592447916666666
$ws.Columns("J").ColumnWidth = 11.166666666666666
$ws.Columns("K").ColumnWidth = 22.451822916666668

# Add a new data row (row 9) with a new SIZE value "2xl".
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = 4
$ws.Cells.Item(9, 3).Value = "JPY"
$ws.Cells.Item(9, 4).Value = "qwerty"
$ws.Cells.Item(9, 6).Value = "blk"
$ws.Cells.Item(9, 7).Value = "2xl"
$ws.Cells.Item(9, 8).Value = 10
$ws.Cells.Item(9, 9).Value = "PCE"
$ws.Cells.Item(9, 10).Value = 5
$ws.Cells.Item(9, 11).Value = 0

# Label the new column last (matches shared-string insertion order: "2xl"
# was registered before "DESCRIPTION").
$ws.Cells.Item(1, 5).Value = "DESCRIPTION"

# Move the active selection like in the saved file.
$ws.Range("D14").Select()
